$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# pandas-based re-export renamed the sheet from the LibreOffice default "Лист1" to "test"
$ws.Name = "test"

# Insert a new column before the old "user" column (E) to host the new "type" column;
# this shifts the former column E ("user" + its "seva" values) to column F.
$ws.Columns("E").Insert()

# Fix the long-standing "weigth" -> "weight" header typo.
$ws.Range("C1").Value = "weight"

# New "type" column header and per-row values.
$ws.Range("E1").Value = "type"
$ws.Range("E2").Value = "S"
$ws.Range("E3").Value = "W"

# The "user" column (now F) gets distinct values instead of the repeated "seva".
$ws.Range("F2").Value = "u1"
$ws.Range("F3").Value = "u2"

# Apply the column widths recorded for the new A:F layout.
$ws.Columns("A").ColumnWidth = 6.43
$ws.Columns("B").ColumnWidth = 7.95
$ws.Columns("C").ColumnWidth = 28.48
$ws.Columns("D").ColumnWidth = 33.76
$ws.Range("E1:F1").EntireColumn.ColumnWidth = 5.04

# Update the active selection to match the saved view state.
$ws.Range("I10").Select()
